$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Hours Burnt" (column G) values for a few tasks (Sprint-2 progress tracking)
$ws.Range("G5").Value = 2
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 1

# Update the view: scroll window back to top, select G11
$ws.Range("G11").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
